$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose leader (D column) flag moves from 1 to 0
$zeroRows = @(2,12,24,32,42,52,64,74,82,92,104,112,122,132,142,152,162,174,182,192,202,214,224,232,244,252,262,272,282,292,302,312,322,332,342,352,364,372,382,392)

# Rows whose leader (D column) flag moves from 0 to 1
$oneRows = @(9,20,31,40,51,59,69,79,91,100,109,119,129,140,151,159,171,180,190,201,209,220,229,239,249,261,270,279,289,299,311,319,330,339,350,360,371,379,389,401)

foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 4).Value = 0
}

foreach ($r in $oneRows) {
    $ws.Cells.Item($r, 4).Value = 1
}

